$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.224.51"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "'1.856.36"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'241.60"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'0.6998"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'0.07779"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.3068"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("D10").Value = "'23.76"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("D11").Value = "'0.07814"
$ws.Range("E11").Value = "  -2.41%  "
$ws.Range("D12").Value = "'1.852.06"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("D13").Value = "'5.103"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").Value = "'92.15"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").Value = "'0.6867"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").Value = "'6.529"
$ws.Range("E16").Value = "  +2.80%  "
$ws.Range("D17").Value = "'0.000008446"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").Value = "'29.212.05"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'250.01"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").Value = "'2.110.35"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "'7.522"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "'0.9999"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "'0.1498"
$ws.Range("E25").Value = "  -3.28%  "
$ws.Range("D26").Value = "'161.05"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("D27").Value = "'8.862"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").Value = "'18.51"
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("D29").Value = "'1.558"
$ws.Range("E29").Value = "  +4.03%  "
$ws.Range("D30").Value = "'4.246"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("D31").Value = "'4.201"
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("D32").Value = "'1.195"
$ws.Range("D33").Value = "'0.05205"
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("E34").Value = "  +2.28%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.168"
$ws.Range("E35").Value = "  +1.17%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'1.842"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").Value = "'0.01861"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "'1.215.08"
$ws.Range("E39").Value = "  -2.89%  "
$ws.Range("D40").Value = "'2.723"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").Value = "'0.8989"
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").Value = "'109.84"
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("D43").Value = "'0.9994"
$ws.Range("D44").Value = "'2.009.19"
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("D45").Value = "'5.505"
$ws.Range("E45").Value = "  -12.04%  "
$ws.Range("D46").Value = "'0.00000000124"
$ws.Range("E46").Value = "  -3.26%  "
$ws.Range("D47").Value = "'65.39"
$ws.Range("E47").Value = "  -8.18%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.558"
$ws.Range("E48").Value = "  +1.80%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.5176"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "'1.754"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("D51").Value = "'7.019"
$ws.Range("E51").Value = "  +0.54%  "
